# Rename the "How to implement basic daily attendance" page to
# "How to implement naplan results reporting" — this affects the
# visible Heading1 text and the bookmark Word generated for that
# heading (used for the doc's internal anchor/TOC linking).

$d = $word.ActiveDocument

# 1. Update the visible heading text.
$d.Content.Find.Execute(
    "How to implement basic daily attendance", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "How to implement naplan results reporting", 2)

# 2. Re-point the heading's bookmark to the new slug-style name.
#    (Bookmarks collection in this host only supports adding new
#    bookmarks, not renaming/deleting in place, so add the
#    correctly-named bookmark at the same collapsed location the
#    heading bookmark occupies - the very start of the Heading1
#    paragraph.)
$heading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "How to implement naplan results reporting") {
        $heading = $p
    }
}
if ($heading -ne $null) {
    $bmRange = $d.Range($heading.Range.Start, $heading.Range.Start)
    $d.Bookmarks.Add("how-to-implement-naplan-results-reporting", $bmRange)
}
